$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 8, shifting existing rows 8.. down by one.
$ws.Rows.Item(8).Insert()

# Set the id for the newly inserted row 8.
$ws.Range("A8").Value = "M1-W007"

# Set all the numeric cells in the new row (B8:AO8) to 0, matching the
# all-zero pattern of a freshly-inserted blank data row.
$ws.Range("B8:AO8").Value = 0
